$d = $word.ActiveDocument
$d.Content.Find.Execute(" WEB APP PARA PSICOLOGO- PACIENTE", $true, $false, $false, $false, $false,
                         $true, 1, $false, " WEB APP PARA PSICOLOGO- PACIENTE.", 2)
